# Generate Report for Handoff
# Updates the "b.md" rows across the Overview / zh-cn / de-de sheets to
# reflect a new handoff of b.md (status changed from "Handed back: in sync
# with en-US" to "Ready for handoff", new handoff file/datetime recorded,
# and an error detail noting the handback file version is stale).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ---------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-12-16 08:31:45"

# ---- zh-cn sheet --------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces Excel to keep "False" as text instead of
# auto-coercing it to a boolean; reset Style afterwards so the cell
# doesn't pick up an extra "quote prefix" number format.
$zh.Range("F3").Value = "'False"
$zh.Range("F3").Style = "Normal"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-12-16 08:31:31"
$zh.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f6bf71ad6fe7a755a3e3be1f95648fff4888ffa/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/94c1eafabe511b213d1e958515ea217be0a63faa/e2e/b.md."
$zh.Range("R1").ColumnWidth = 39.14285714285714

# ---- de-de sheet --------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "'False"
$de.Range("F3").Style = "Normal"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-12-16 08:31:45"
$de.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f6bf71ad6fe7a755a3e3be1f95648fff4888ffa/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/94c1eafabe511b213d1e958515ea217be0a63faa/e2e/b.md."
$de.Range("R1").ColumnWidth = 39.14285714285714
